$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen column A (account for runtime's internal width-unit padding offset
# so the stored OOXML width ends up as exactly 80)
$ws.Columns.Item(1).ColumnWidth = 79.16666666666667

# New Q&A rows appended after row 40
$qa = @(
    @("How many different lithology types can be recorded in a log at most?", "The highest number of lithology types that can be represented in a log is 450."),
    @("What is the upper limit on lithology types in a log?", "The maximum number of lithology types that can be recorded in a log is 450."),
    @("What's the highest number of lithology types that can be represented in a log?", "The maximum number of lithology types that can be recorded in a log is 450."),
    @("How many lithology types can a single log accommodate at maximum?", "The maximum number of lithology types that can be recorded in a log is 450.")
)

$row = 41
foreach ($pair in $qa) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
